$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Doctors")

# Row 2
$ws.Range("A2").Value = "Dr. S Sarangapani"
$ws.Range("B2").Value = "Ayurveda"
$ws.Range("C2").Value = "44 years experience overall"
$ws.Range("D2").Value = "Erragadda,Hyderabad"

# Row 3
$ws.Range("A3").Value = "Dr. M. Narasimha"
$ws.Range("B3").Value = "Ayurveda"
$ws.Range("C3").Value = "25 years experience overall"
$ws.Range("D3").Value = "Ameerpet,Hyderabad"

# Row 4
$ws.Range("A4").Value = "Dr. C Narmada"
$ws.Range("B4").Value = "Ayurveda"
$ws.Range("C4").Value = "24 years experience overall"
$ws.Range("D4").Value = "Shaikpet,Hyderabad"

# Row 5
$ws.Range("A5").Value = "Dr. Priti Thakre"
$ws.Range("B5").Value = "Ayurveda"
$ws.Range("C5").Value = "22 years experience overall"
$ws.Range("D5").Value = "KPHB,Hyderabad"

# Row 6
$ws.Range("A6").Value = "Dr. Sunita Grace"
$ws.Range("B6").Value = "Ayurveda"
$ws.Range("C6").Value = "21 years experience overall"
$ws.Range("D6").Value = "Banjara Hills,Hyderabad"
